# Weekly update: insert a new price-report row for the week of 2023-07-25
# (date serial 45132) just before the existing row 525, shifting the rest
# of the "Zanahoria" (carrot) records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 525; everything currently at/after row 525
# (through the old last row 572) shifts down to 526..573.
$ws.Rows.Item(525).EntireRow.Insert()

# Populate the newly inserted row 525 with this week's data.
$ws.Cells.Item(525, 1).Value  = 8
$ws.Cells.Item(525, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(525, 3).Value  = "Coquimbo"
$ws.Cells.Item(525, 4).Value  = 45132
$ws.Cells.Item(525, 5).Value  = 4
$ws.Cells.Item(525, 6).Value  = 100114013
$ws.Cells.Item(525, 7).Value  = "Zanahoria"
$ws.Cells.Item(525, 8).Value  = "Sin especificar"
$ws.Cells.Item(525, 9).Value  = "Primera"
$ws.Cells.Item(525, 10).Value = 500
$ws.Cells.Item(525, 11).Value = 5800
$ws.Cells.Item(525, 12).Value = 6000
$ws.Cells.Item(525, 13).Value = 5900
$ws.Cells.Item(525, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(525, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(525, 16).Value = 295
$ws.Cells.Item(525, 17).Value = 20
$ws.Cells.Item(525, 18).Value = "Hortaliza"
